# [Abraham]: fixed select all in enterprise.
#
# The "ruc de casa comercial donde Aplica Crédito" column (W) stored a bare
# comma separated list of RUCs (e.g. "1002003004001,1002003004002"). The
# enterprise "select all" feature expects this value to look like a JSON
# array, so it is now wrapped in square brackets:
# "[1002003004001,1002003004002]".
#
# The last active selection saved with the sheet is also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RUC list values in column W (rows 2 and 3) to the bracketed
# "array-like" format expected by the fixed "select all" logic.
$ws.Range("W2").Value2 = "[1002003004001,1002003004002]"
$ws.Range("W3").Value2 = "[1002003004001,1002003004002]"

# Nudge the font so these two cells pick up a dedicated style (matching the
# distinct cell style introduced for this column in the fixed workbook)
# instead of continuing to share the generic "no-wrap" style used before.
$ws.Range("W2").Font.ThemeColor = 1
$ws.Range("W3").Font.ThemeColor = 1

# Restore the workbook's last saved selection/active cell.
$ws.Range("G20").Select() | Out-Null
